$wb = $excel.ActiveWorkbook

$newTimestamp = "05/06/2024 17:10:12"

# ---------------------------------------------------------------------------
# Sheet "10per change": 9 data rows (2..10) -> duplicate as rows 11..19,
# then convert the bsecode column (D) of the original rows to real numbers.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("10per change")

$srcRows1 = 9
$firstDataRow1 = 2
$lastDataRow1 = $firstDataRow1 + $srcRows1 - 1          # 10
$firstNewRow1 = $lastDataRow1 + 1                        # 11
$lastNewRow1 = $firstNewRow1 + $srcRows1 - 1              # 19

$srcRange1 = "A" + $firstDataRow1 + ":H" + $lastDataRow1
$dstRange1 = "A" + $firstNewRow1 + ":H" + $lastNewRow1

$ws1.Range($srcRange1).Copy() | Out-Null
$ws1.Range($dstRange1).PasteSpecial() | Out-Null

# Stamp the freshly duplicated rows with the later scrape time.
for ($r = $firstNewRow1; $r -le $lastNewRow1; $r++) {
    $ws1.Cells.Item($r, 8).Value = $newTimestamp
}

# Convert the original rows' bsecode (column D) from text to numeric.
for ($r = $firstDataRow1; $r -le $lastDataRow1; $r++) {
    $cell = $ws1.Cells.Item($r, 4)
    $cell.Value = $cell.Value2
}

# ---------------------------------------------------------------------------
# Sheet "DND 3 V 0.3": 1 data row (2) -> duplicate as row 3, then convert the
# bsecode column (D) of the original row to a real number.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("DND 3 V 0.3")

$firstDataRow3 = 2
$lastDataRow3 = 2
$firstNewRow3 = 3
$lastNewRow3 = 3

$srcRange3 = "A" + $firstDataRow3 + ":H" + $lastDataRow3
$dstRange3 = "A" + $firstNewRow3 + ":H" + $lastNewRow3

$ws3.Range($srcRange3).Copy() | Out-Null
$ws3.Range($dstRange3).PasteSpecial() | Out-Null

for ($r = $firstNewRow3; $r -le $lastNewRow3; $r++) {
    $ws3.Cells.Item($r, 8).Value = $newTimestamp
}

for ($r = $firstDataRow3; $r -le $lastDataRow3; $r++) {
    $cell = $ws3.Cells.Item($r, 4)
    $cell.Value = $cell.Value2
}
